$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 54
$ws.Range("C2").Value = "dog/dog025.jpg"
$ws.Range("D2").Value = "fesseln"
$ws.Range("E2").Value = "dog"
$ws.Range("B3").Value = 123
$ws.Range("C3").Value = "dog/dog001.jpg"
$ws.Range("D3").Value = "liefern"
$ws.Range("E3").Value = "dog"
$ws.Range("B4").Value = 62
$ws.Range("C4").Value = "dog/dog024.jpg"
$ws.Range("D4").Value = "saufen"
$ws.Range("E4").Value = "dog"
$ws.Range("B5").Value = 79
$ws.Range("C5").Value = "dog/dog023.jpg"
$ws.Range("D5").Value = "hoffen"
$ws.Range("E5").Value = "dog"
$ws.Range("B6").Value = 81
$ws.Range("C6").Value = "flower/flower016.jpg"
$ws.Range("D6").Value = "schätzen"
$ws.Range("E6").Value = "flower"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = "dog/dog014.jpg"
$ws.Range("D7").Value = "sieben"
$ws.Range("E7").Value = "dog"
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = "dog/dog027.jpg"
$ws.Range("D8").Value = "enden"
$ws.Range("E8").Value = "dog"
$ws.Range("B9").Value = 71
$ws.Range("C9").Value = "dog/dog020.jpg"
$ws.Range("D9").Value = "drehen"
$ws.Range("E9").Value = "dog"
$ws.Range("B10").Value = 121
$ws.Range("C10").Value = "dog/dog000.jpg"
$ws.Range("D10").Value = "dauern"
$ws.Range("E10").Value = "dog"
$ws.Range("B11").Value = 14
$ws.Range("C11").Value = "dog/dog009.jpg"
$ws.Range("D11").Value = "füllen"
$ws.Range("E11").Value = "dog"
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = "dog/dog013.jpg"
$ws.Range("D12").Value = "fliehen"
$ws.Range("E12").Value = "dog"
$ws.Range("B13").Value = 117
$ws.Range("C13").Value = "dog/dog029.jpg"
$ws.Range("D13").Value = "haken"
$ws.Range("E13").Value = "dog"
$ws.Range("B14").Value = 39
$ws.Range("C14").Value = "flower/flower023.jpg"
$ws.Range("D14").Value = "lehnen"
$ws.Range("E14").Value = "flower"
$ws.Range("B15").Value = 75
$ws.Range("C15").Value = "flower/flower026.jpg"
$ws.Range("D15").Value = "hupen"
$ws.Range("E15").Value = "flower"
$ws.Range("B16").Value = 67
$ws.Range("C16").Value = "dog/dog015.jpg"
$ws.Range("D16").Value = "gründen"
$ws.Range("E16").Value = "dog"
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = "flower/flower022.jpg"
$ws.Range("D17").Value = "opfern"
$ws.Range("E17").Value = "flower"
$ws.Range("B18").Value = 16
$ws.Range("C18").Value = "flower/flower025.jpg"
$ws.Range("D18").Value = "strahlen"
$ws.Range("E18").Value = "flower"
$ws.Range("B19").Value = 60
$ws.Range("C19").Value = "flower/flower006.jpg"
$ws.Range("D19").Value = "tagen"
$ws.Range("E19").Value = "flower"
$ws.Range("B20").Value = 87
$ws.Range("C20").Value = "flower/flower005.jpg"
$ws.Range("D20").Value = "kehren"
$ws.Range("E20").Value = "flower"
$ws.Range("B21").Value = 64
$ws.Range("C21").Value = "flower/flower021.jpg"
$ws.Range("D21").Value = "stärken"
$ws.Range("E21").Value = "flower"
$ws.Range("B22").Value = 107
$ws.Range("C22").Value = "flower/flower000.jpg"
$ws.Range("D22").Value = "backen"
$ws.Range("E22").Value = "flower"
$ws.Range("B23").Value = 109
$ws.Range("C23").Value = "flower/flower001.jpg"
$ws.Range("D23").Value = "schicken"
$ws.Range("E23").Value = "flower"
$ws.Range("B24").Value = 59
$ws.Range("C24").Value = "flower/flower014.jpg"
$ws.Range("D24").Value = "tauschen"
$ws.Range("E24").Value = "flower"
$ws.Range("B25").Value = 118
$ws.Range("C25").Value = "dog/dog012.jpg"
$ws.Range("D25").Value = "drohen"
$ws.Range("E25").Value = "dog"
$ws.Range("B26").Value = 48
$ws.Range("C26").Value = "flower/flower011.jpg"
$ws.Range("D26").Value = "raten"
$ws.Range("E26").Value = "flower"
$ws.Range("B27").Value = 7
$ws.Range("C27").Value = "dog/dog007.jpg"
$ws.Range("D27").Value = "sondern"
$ws.Range("E27").Value = "dog"
$ws.Range("B28").Value = 52
$ws.Range("C28").Value = "flower/flower031.jpg"
$ws.Range("D28").Value = "biegen"
$ws.Range("E28").Value = "flower"
$ws.Range("B29").Value = 11
$ws.Range("C29").Value = "dog/dog022.jpg"
$ws.Range("D29").Value = "hauen"
$ws.Range("E29").Value = "dog"
$ws.Range("B30").Value = 26
$ws.Range("C30").Value = "flower/flower004.jpg"
$ws.Range("D30").Value = "bitten"
$ws.Range("E30").Value = "flower"
$ws.Range("B31").Value = 91
$ws.Range("C31").Value = "flower/flower013.jpg"
$ws.Range("D31").Value = "posten"
$ws.Range("E31").Value = "flower"
$ws.Range("B32").Value = 12
$ws.Range("C32").Value = "dog/dog031.jpg"
$ws.Range("D32").Value = "töten"
$ws.Range("E32").Value = "dog"
$ws.Range("B33").Value = 9
$ws.Range("C33").Value = "flower/flower003.jpg"
$ws.Range("D33").Value = "ehren"
$ws.Range("E33").Value = "flower"
